$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cell H1, copying style from existing header (G1) so it matches
# the bold/centered/bordered look of the other header cells.
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122) # xlPasteFormats

# Add data values in the new column for the two data rows
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 1
